# Hydraulic gear system additions (#6893)
# Adds "Nose gear door" and "Nose gear" actuator data blocks to the
# "Actuators" sheet, adds max-force helper columns (H/I) to the existing
# actuator blocks, tweaks the "Main gear door" actuator bore/rod sizing,
# adds threaded review comments, and repositions the landing-gear-door
# clearance illustration.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actuators")

# ---------------------------------------------------------------------
# 1) Add "Max force retract" / "Max force extend" helper columns to the
#    existing actuator blocks (Gear Main, the door block at 36-40, and
#    Main gear door at 43-47).
# ---------------------------------------------------------------------

# "Gear Main" block (row 30 data row) -- note: labels intentionally
# mirror the source workbook's own (swapped) column order.
$ws.Range("H32").Formula = "=""Max force extend"""
$ws.Range("I32").Formula = "=""Max force retract"""
$ws.Range("H33").Formula = "=(E30-F30)*20684000"
$ws.Range("I33").Formula = "=E30*20684000"

# Door block (row 37 data row)
$ws.Range("H39").Formula = "=""Max force retract"""
$ws.Range("I39").Formula = "=""Max force extend"""
$ws.Range("H40").Formula = "=(E37-F37)*20684000"
$ws.Range("I40").Formula = "=E37*20684000"

# ---------------------------------------------------------------------
# 2) Main gear door actuator (row 44): corrected bore / rod sizing.
# ---------------------------------------------------------------------
$ws.Range("C44").Value = 0.0538
$ws.Range("D44").Formula = "=0.03015"

$ws.Range("H46").Formula = "=""Max force retract"""
$ws.Range("I46").Formula = "=""Max force extend"""
$ws.Range("H47").Formula = "=(E44-F44)*20684000"
$ws.Range("I47").Formula = "=E44*20684000"

# ---------------------------------------------------------------------
# 3) New "Nose gear door" actuator block (rows 50-56), modeled after the
#    existing "Main gear door" block (rows 43-47).
# ---------------------------------------------------------------------
$ws.Range("A43:J47").Copy()
$ws.Range("A50:J54").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B50").Formula = "=B43"
$ws.Range("C50").Formula = "=C43"
$ws.Range("D50").Formula = "=D43"
$ws.Range("E50").Formula = "=E43"
$ws.Range("F50").Formula = "=F43"
$ws.Range("G50").Formula = "=G43"
$ws.Range("H50").Formula = "=H43"
$ws.Range("I50").Formula = "=I43"
$ws.Range("J50").Formula = "=J43"

$ws.Range("A51").Value = "Nose gear door"
$ws.Range("B51").Value = 0.162
$ws.Range("C51").Value = 0.0378
$ws.Range("D51").Formula = "=0.023"
$ws.Range("E51").Formula = "=PI()* (C51/2)^2"
$ws.Range("F51").Formula = "=PI()* (D51/2)^2"
$ws.Range("G51").Formula = "=E51*B51"
$ws.Range("H51").Formula = "=(E51-F51)*B51"
$ws.Range("I51").Formula = "=G51/H51"
$ws.Range("J51").Formula = "=(G51-H51) * 1000"

$ws.Range("A52").Formula = "=A45"
$ws.Range("G52").Formula = "=G51*2"
$ws.Range("H52").Formula = "=H51*2"
$ws.Range("I52").Formula = "=G52/H52"
$ws.Range("J52").Formula = "=(G52-H52) * 1000"

$ws.Range("B53").Formula = "=B46"
$ws.Range("C53").Formula = "=C46"
$ws.Range("D53").Formula = "=D46"
$ws.Range("E53").Formula = "=E46"
$ws.Range("F53").Formula = "=F46"
$ws.Range("H53").Formula = "=H46"
$ws.Range("I53").Formula = "=I46"

$ws.Range("C54").Value = -0.1465
$ws.Range("D54").Value = 0
$ws.Range("E54").Value = -0.1465
$ws.Range("F54").Value = 0.4
$ws.Range("H54").Formula = "=(E51-F51)*20684000"
$ws.Range("I54").Formula = "=E51*20684000"

$ws.Range("C55").Value = "Note those coordinates are chosen ""randomly"" so we get the correct actuator travel of 162mm"
$ws.Range("C56").Value = "This actuator in fact should push a crankbell that gives a particular motion ratio."

# ---------------------------------------------------------------------
# 4) New "Nose gear" actuator block (rows 58-63).
# ---------------------------------------------------------------------
$ws.Range("A43:J47").Copy()
$ws.Range("A58:J62").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B58").Formula = "=B43"
$ws.Range("C58").Formula = "=C43"
$ws.Range("D58").Formula = "=D43"
$ws.Range("E58").Formula = "=E43"
$ws.Range("F58").Formula = "=F43"
$ws.Range("G58").Formula = "=G43"
$ws.Range("H58").Formula = "=H43"
$ws.Range("I58").Formula = "=I43"
$ws.Range("J58").Formula = "=J43"

$ws.Range("A59").Value = "Nose gear"
$ws.Range("B59").Value = 0.32
$ws.Range("C59").Value = 0.0792
$ws.Range("D59").Formula = "=0.035"
$ws.Range("E59").Formula = "=PI()* (C59/2)^2"
$ws.Range("F59").Formula = "=PI()* (D59/2)^2"
$ws.Range("G59").Formula = "=E59*B59"
$ws.Range("H59").Formula = "=(E59-F59)*B59"
$ws.Range("I59").Formula = "=G59/H59"
$ws.Range("J59").Formula = "=(G59-H59) * 1000"

$ws.Range("A60").Formula = "=A45"
$ws.Range("G60").Formula = "=G59*2"
$ws.Range("H60").Formula = "=H59*2"
$ws.Range("I60").Formula = "=G60/H60"
$ws.Range("J60").Formula = "=(G60-H60) * 1000"

$ws.Range("B61").Formula = "=B46"
$ws.Range("C61").Value = "Control Arm position Z"
$ws.Range("D61").Formula = "=D46"
$ws.Range("E61").Value = "Anchor point position Z"
$ws.Range("F61").Formula = "=F46"
$ws.Range("H61").Formula = "=H46"
$ws.Range("I61").Formula = "=I46"

$ws.Range("C62").Value = 0.212
$ws.Range("D62").Value = -0.093
$ws.Range("E62").Value = 0
$ws.Range("F62").Value = 0.56
$ws.Range("H62").Formula = "=(E59-F59)*20684000"
$ws.Range("I62").Formula = "=E59*20684000"

$ws.Range("C63").Value = "Note those coordinates are chosen ""randomly"" so we get the correct actuator travel of 320mm"

# ---------------------------------------------------------------------
# 5) Threaded review comments on the two new actuator values.
# ---------------------------------------------------------------------
$ws.Range("D54").AddCommentThreaded("This value gives correct actuator displacement")
$ws.Range("D62").AddCommentThreaded("This value gives correct actuator displacement")

# ---------------------------------------------------------------------
# 6) Reposition the landing-gear-door clearance illustration further
#    down/right on the sheet so it sits beside the new data blocks.
# ---------------------------------------------------------------------
$pic = $ws.Shapes.Item("Image 3")
$pic.Left = 791.88
$pic.Top = 402.86
$pic.LockAspectRatio = $false
$pic.Width = 283.2
$pic.Height = 306.5

$wb.Save()
